# Update the cryptos worksheet with refreshed price / 1h-volume figures.
# Column D holds prices (stored as text, some of which look like plain
# numbers) and column E holds the "  +/-X.XX%  " volume-change text.
#
# For D-column values that parse as a plain number (e.g. "2.90", "6.00",
# "0.712"), a leading apostrophe is used so Excel keeps them as literal
# text (preserving trailing zeros / exact formatting) instead of coercing
# them into numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.334.96'
$ws.Range("E2").Value = '  -3.56%  '

$ws.Range("D3").Value = '3.164.76'
$ws.Range("E3").Value = '  -2.82%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'607.46"
$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").Value = "'147.87"
$ws.Range("E6").Value = '  -6.62%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '3.158.02'
$ws.Range("E8").Value = '  -3.01%  '

$ws.Range("E9").Value = '  -4.08%  '

$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = '  -6.50%  '

$ws.Range("D11").Value = "'5.53"
$ws.Range("E11").Value = '  -7.13%  '

$ws.Range("E12").Value = '  -6.15%  '

$ws.Range("E13").Value = '  -7.13%  '

$ws.Range("D14").Value = "'35.74"
$ws.Range("E14").Value = '  -9.50%  '

$ws.Range("D15").Value = '3.684.89'
$ws.Range("E15").Value = '  -2.80%  '

$ws.Range("D16").Value = '64.328.99'
$ws.Range("E16").Value = '  -3.58%  '

$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").Value = '3.165.33'
$ws.Range("E18").Value = '  -2.64%  '

$ws.Range("D19").Value = "'6.96"
$ws.Range("E19").Value = '  -6.16%  '

$ws.Range("D20").Value = "'482.45"
$ws.Range("E20").Value = '  -5.05%  '

$ws.Range("D21").Value = "'14.79"
$ws.Range("E21").Value = '  -3.92%  '

$ws.Range("D22").Value = "'0.712"
$ws.Range("E22").Value = '  -5.50%  '

$ws.Range("D23").Value = "'7.76"
$ws.Range("E23").Value = '  -4.10%  '

$ws.Range("E24").Value = '  -7.30%  '

$ws.Range("D25").Value = "'83.70"
$ws.Range("E25").Value = '  -3.31%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("E27").Value = '  -4.98%  '

$ws.Range("D28").Value = "'8.50"
$ws.Range("E28").Value = '  -6.46%  '

$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = '  -7.98%  '

$ws.Range("D30").Value = "'6.83"
$ws.Range("E30").Value = '  -0.71%  '

$ws.Range("E31").Value = '  -23.22%  '

$ws.Range("E32").Value = '  -3.30%  '

$ws.Range("E33").Value = '  +0.15%  '

$ws.Range("D34").Value = "'26.29"
$ws.Range("E34").Value = '  -6.52%  '

$ws.Range("E35").Value = '  -4.69%  '

$ws.Range("D36").Value = "'6.00"
$ws.Range("E36").Value = '  -6.81%  '

$ws.Range("E37").Value = '  -2.09%  '

$ws.Range("E38").Value = '  -9.39%  '

$ws.Range("D39").Value = "'452.87"
$ws.Range("E39").Value = '  -8.53%  '

$ws.Range("D40").Value = "'2.90"
$ws.Range("E40").Value = '  -13.53%  '

$ws.Range("D42").Value = "'8.46"
$ws.Range("E42").Value = '  -4.30%  '

$ws.Range("E43").Value = '  -7.90%  '

$ws.Range("D44").Value = '2.854.40'
$ws.Range("E44").Value = '  -3.30%  '

$ws.Range("E45").Value = '  -9.38%  '

$ws.Range("E46").Value = '  -8.62%  '

$ws.Range("D47").Value = "'26.45"
$ws.Range("E47").Value = '  -7.65%  '

$ws.Range("E49").Value = '  -6.84%  '

$ws.Range("E50").Value = '  -4.55%  '

$ws.Range("D51").Value = "'119.32"
$ws.Range("E51").Value = '  -1.96%  '
